$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 5 ("Duy ", 44650, 5) is being pushed down to row 9 to make
# room for four newly logged entries (rows 5-8). Write the relocated row
# first so the "Duy " shared string stays referenced while we overwrite row 5.
$ws.Range("A9").Value = "Duy "
$ws.Range("B9").NumberFormat = "d\-mmm"
$ws.Range("B9").Value = 44650
$ws.Range("C9").Value = 5

# New entries logged for Eric (and one more "All" entry).
$ws.Range("A5").Value = "Eric"
$ws.Range("B5").NumberFormat = "d\-mmm"
$ws.Range("B5").Value = 44647
$ws.Range("C5").Value = 2

$ws.Range("A6").Value = "All"
$ws.Range("B6").NumberFormat = "d\-mmm"
$ws.Range("B6").Value = 44648
$ws.Range("C6").Value = 1

$ws.Range("A7").Value = "Eric "
$ws.Range("B7").NumberFormat = "d\-mmm"
$ws.Range("B7").Value = 44649
$ws.Range("C7").Value = 1

$ws.Range("A8").Value = "Eric"
$ws.Range("B8").NumberFormat = "d\-mmm"
$ws.Range("B8").Value = 44650
$ws.Range("C8").Value = 3

# Extend the styled (but otherwise empty) tracking rows at the bottom by two.
$ws.Range("B19").NumberFormat = "d\-mmm"
$ws.Range("B20").NumberFormat = "d\-mmm"

# Selection moved to D11 in the saved file.
$ws.Range("D11").Select()
